$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $escaped = $val -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue 'D2' '29.553.51'
Set-TextValue 'E2' '  +0.07%  '
Set-TextValue 'D3' '1.913.68'
Set-TextValue 'E3' '  -0.21%  '
Set-TextValue 'E4' '  +0.77%  '
Set-TextValue 'D5' '325.64'
Set-TextValue 'E5' '  -0.30%  '
Set-TextValue 'E6' '  +0.65%  '
Set-TextValue 'D7' '0.4825'
Set-TextValue 'E7' '  +0.68%  '
Set-TextValue 'D8' '0.4071'
Set-TextValue 'E8' '  -0.80%  '
Set-TextValue 'D9' '0.08147'
Set-TextValue 'E9' '  +1.12%  '
Set-TextValue 'D10' '1.012'
Set-TextValue 'E10' '  +0.04%  '
Set-TextValue 'D11' '23.42'
Set-TextValue 'E11' '  +4.44%  '
Set-TextValue 'D12' '1.930.91'
Set-TextValue 'E12' '  +1.56%  '
Set-TextValue 'D13' '6.017'
Set-TextValue 'E13' '  +1.05%  '
Set-TextValue 'D14' '7.103'
Set-TextValue 'E14' '  -0.76%  '
Set-TextValue 'D15' '90.38'
Set-TextValue 'E15' '  +0.81%  '
Set-TextValue 'B16' 'BinanceUSD'
Set-TextValue 'C16' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D16' '1.008'
Set-TextValue 'E16' '  +0.73%  '
Set-TextValue 'B17' 'TRON'
Set-TextValue 'C17' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D17' '0.06771'
Set-TextValue 'E17' '  +2.37%  '
Set-TextValue 'D18' '0.00001040'
Set-TextValue 'E18' '  +0.72%  '
Set-TextValue 'D19' '17.71'
Set-TextValue 'E19' '  -0.29%  '
Set-TextValue 'E20' '  +0.63%  '
Set-TextValue 'D21' '29.568.84'
Set-TextValue 'E21' '  +0.08%  '
Set-TextValue 'D22' '5.619'
Set-TextValue 'E22' '  +1.22%  '
Set-TextValue 'D23' '11.82'
Set-TextValue 'E23' '  +2.62%  '
Set-TextValue 'E24' '  -1.33%  '
Set-TextValue 'D25' '2.132.90'
Set-TextValue 'E25' '  +0.07%  '
Set-TextValue 'D26' '154.77'
Set-TextValue 'E26' '  +0.99%  '
Set-TextValue 'D27' '20.04'
Set-TextValue 'E27' '  +0.91%  '
Set-TextValue 'D28' '6.299'
Set-TextValue 'E28' '  +8.26%  '
Set-TextValue 'D29' '2.101'
Set-TextValue 'E29' '  -1.91%  '
Set-TextValue 'D30' '119.79'
Set-TextValue 'E30' '  +1.55%  '
Set-TextValue 'E31' '  -2.87%  '
Set-TextValue 'D32' '0.09554'
Set-TextValue 'E32' '  +0.09%  '
Set-TextValue 'D33' '5.518'
Set-TextValue 'E33' '  +2.38%  '
Set-TextValue 'E34' '  -0.23%  '
Set-TextValue 'D35' '1.392'
Set-TextValue 'E35' '  -2.50%  '
Set-TextValue 'D36' '0.02266'
Set-TextValue 'E36' '  +0.38%  '
Set-TextValue 'D37' '0.06105'
Set-TextValue 'E37' '  +0.02%  '
Set-TextValue 'D38' '1.174'
Set-TextValue 'E38' '  -0.21%  '
Set-TextValue 'D39' '10.82'
Set-TextValue 'E39' '  +6.47%  '
Set-TextValue 'D40' '0.5936'
Set-TextValue 'E40' '  +0.78%  '
Set-TextValue 'D41' '7.923'
Set-TextValue 'E41' '  -5.01%  '
Set-TextValue 'D42' '0.1854'
Set-TextValue 'E42' '  +0.30%  '
Set-TextValue 'D43' '2.455'
Set-TextValue 'E43' '  -2.80%  '
Set-TextValue 'D44' '1.284'
Set-TextValue 'E44' '  -0.67%  '
Set-TextValue 'D45' '0.07732'
Set-TextValue 'E45' '  -3.41%  '
Set-TextValue 'D46' '12.42'
Set-TextValue 'E46' '  +2.41%  '
Set-TextValue 'D47' '0.5566'
Set-TextValue 'E47' '  +0.28%  '
Set-TextValue 'D48' '1.947'
Set-TextValue 'E48' '  +0.68%  '
Set-TextValue 'D49' '115.76'
Set-TextValue 'E49' '  +2.05%  '
Set-TextValue 'D50' '72.72'
Set-TextValue 'E50' '  +1.40%  '
Set-TextValue 'D51' '1.053'
Set-TextValue 'E51' '  +1.72%  '

$excel.CutCopyMode = $false

